$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.696.37"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.616.17"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.95"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.07"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.112"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.80"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.390"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.94"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "3.087.87"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "63.505.57"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000162"
$ws.Range("E16").Value = "  +9.01%  "
$ws.Range("D17").Value = "2.614.97"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.18"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.78"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.68"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.95"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.15"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.70"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.25"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.68"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "553.52"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.163"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "0.0₃0872"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.57"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.13"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.03"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.417"
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.59"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.03"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.62"
$ws.Range("E44").Value = "  +10.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0583"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  +8.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.636"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0252"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0963"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.33"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +18.70%  "
